# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted into the "Berenjena" sheet.
# It lands right after the existing row 108 (i.e. becomes the new row 109),
# pushing every subsequent record down by one row (old row 109 -> 110,
# old row 110 -> 111, ... old row 136 -> 137). The sheet's used range grows
# from A1:R136 to A1:R137 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: insert a blank row at 109, shifting the
# existing rows 109-136 down to 110-137.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(109, 1).Value = 8
$ws.Cells.Item(109, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(109, 3).Value = "Coquimbo"
$ws.Cells.Item(109, 4).Value = 44711
$ws.Cells.Item(109, 5).Value = 4
$ws.Cells.Item(109, 6).Value = 100112001
$ws.Cells.Item(109, 7).Value = "Berenjena"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 600
$ws.Cells.Item(109, 11).Value = 8000
$ws.Cells.Item(109, 12).Value = 9000
$ws.Cells.Item(109, 13).Value = 8500
$ws.Cells.Item(109, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(109, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(109, 16).Value = 170
$ws.Cells.Item(109, 17).Value = 50
$ws.Cells.Item(109, 18).Value = "Hortaliza"
